$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 631
$ws.Range("G2").Value = 1258
$ws.Range("H2").Value = 348
$ws.Range("J2").Value = 11
$ws.Range("N2").Value = 27215287.94
$ws.Range("O2").Value = 45018323.74
$ws.Range("P2").Value = 11540812.19
$ws.Range("R2").Value = 505017.99
$ws.Range("V2").Value = 32.29
$ws.Range("W2").Value = 53.42
$ws.Range("X2").Value = 13.69
$ws.Range("AB2").Value = 0.6
$ws.Range("AD2").Value = 7054848.94
$ws.Range("AE2").Value = 67.11
$ws.Range("AF2").Value = 13.69

# Row 3
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1
$ws.Range("N3").Value = 57565.2
$ws.Range("O3").Value = 56835.24
$ws.Range("T3").Value = 10366.72
$ws.Range("V3").Value = 46.14
$ws.Range("W3").Value = 45.55
$ws.Range("Z3").Value = 8.31
$ws.Range("AD3").Value = 16052
$ws.Range("AE3").Value = 53.86
$ws.Range("AF3").Value = 8.31

# Row 4
$ws.Range("F4").Value = 1
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = 16852.6
$ws.Range("T4").Value = 0
$ws.Range("V4").Value = 100
$ws.Range("Z4").Value = 0
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0

# Row 5
$ws.Range("F5").Value = 373
$ws.Range("G5").Value = 232
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 37
$ws.Range("N5").Value = 12287505.76
$ws.Range("O5").Value = 7357918.92
$ws.Range("P5").Value = 126342.74
$ws.Range("Q5").Value = 136181.6
$ws.Range("R5").Value = 0
$ws.Range("T5").Value = 547957.09
$ws.Range("V5").Value = 60.07
$ws.Range("W5").Value = 35.97
$ws.Range("X5").Value = 0.62
$ws.Range("Z5").Value = 2.68
$ws.Range("AA5").Value = 0.67
$ws.Range("AB5").Value = 0
$ws.Range("AD5").Value = 1261361.67
$ws.Range("AE5").Value = 39.94
$ws.Range("AF5").Value = 3.97
